# Update attendee/view counts (column F) on sheets "展览", "演出", and "全部类型"
# per the upstream gh-pages data refresh (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 188
$ws.Range("F5").Value  = 1093
$ws.Range("F6").Value  = 8221
$ws.Range("F7").Value  = 8221
$ws.Range("F8").Value  = 140
$ws.Range("F9").Value  = 211
$ws.Range("F10").Value = 6906
$ws.Range("F12").Value = 5025
$ws.Range("F13").Value = 5480
$ws.Range("F15").Value = 334
$ws.Range("F25").Value = 9232
$ws.Range("F27").Value = 1683
$ws.Range("F28").Value = 877
$ws.Range("F31").Value = 1876
$ws.Range("F37").Value = 1893
$ws.Range("F39").Value = 1205
$ws.Range("F41").Value = 4820
$ws.Range("F44").Value = 78
$ws.Range("F45").Value = 4
$ws.Range("F49").Value = 922
$ws.Range("F50").Value = 1265
$ws.Range("F51").Value = 43

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F17").Value = 894

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 188
$ws.Range("F6").Value  = 1094
$ws.Range("F7").Value  = 8222
$ws.Range("F8").Value  = 140
$ws.Range("F9").Value  = 211
$ws.Range("F10").Value = 6906
$ws.Range("F14").Value = 5025
$ws.Range("F15").Value = 5480
$ws.Range("F17").Value = 334
$ws.Range("F25").Value = 9232
$ws.Range("F27").Value = 1683
$ws.Range("F28").Value = 877
$ws.Range("F31").Value = 1876
$ws.Range("F37").Value = 1893
$ws.Range("F39").Value = 1205
$ws.Range("F41").Value = 4820
$ws.Range("F44").Value = 78
$ws.Range("F45").Value = 4
$ws.Range("F49").Value = 922
$ws.Range("F50").Value = 1265
$ws.Range("F51").Value = 43
